$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(150,1).NumberFormat = "@"
$ws.Cells.Item(150,1).Value = "02-08-2021"
$ws.Cells.Item(150,2).Value = 2.57
$ws.Cells.Item(150,3).Value = 2.91
$ws.Cells.Item(150,4).Value = 3.19
$ws.Cells.Item(150,5).Value = 3.44
$ws.Cells.Item(150,6).Value = 4.22
$ws.Cells.Item(150,7).Value = -1.36
$ws.Cells.Item(150,8).Value = -0.77
$ws.Cells.Item(150,9).Value = -0.37
$ws.Cells.Item(150,10).Value = -0.03
$ws.Cells.Item(150,11).Value = 0.22
$ws.Cells.Item(150,12).Value = 1.02
$ws.Cells.Item(150,13).Value = 1.25

$ws.Cells.Item(151,1).NumberFormat = "@"
$ws.Cells.Item(151,1).Value = "03-08-2021"
$ws.Cells.Item(151,2).Value = 2.63
$ws.Cells.Item(151,3).Value = 2.98
$ws.Cells.Item(151,4).Value = 3.27
$ws.Cells.Item(151,5).Value = 3.52
$ws.Cells.Item(151,6).Value = 4.28
$ws.Cells.Item(151,7).Value = -1.36
$ws.Cells.Item(151,8).Value = -0.73
$ws.Cells.Item(151,9).Value = -0.32
$ws.Cells.Item(151,10).Value = 0.01
$ws.Cells.Item(151,11).Value = 0.27
$ws.Cells.Item(151,12).Value = 1.08
$ws.Cells.Item(151,13).Value = 1.3

$ws.Cells.Item(152,1).NumberFormat = "@"
$ws.Cells.Item(152,1).Value = "04-08-2021"
$ws.Cells.Item(152,2).Value = 2.66
$ws.Cells.Item(152,3).Value = 3.01
$ws.Cells.Item(152,4).Value = 3.31
$ws.Cells.Item(152,5).Value = 3.56
$ws.Cells.Item(152,6).Value = 4.34
$ws.Cells.Item(152,7).Value = -1.34
$ws.Cells.Item(152,8).Value = -0.7
$ws.Cells.Item(152,9).Value = -0.28
$ws.Cells.Item(152,10).Value = 0.05
$ws.Cells.Item(152,11).Value = 0.3
$ws.Cells.Item(152,12).Value = 1.11
$ws.Cells.Item(152,13).Value = 1.33

$ws.Cells.Item(153,1).NumberFormat = "@"
$ws.Cells.Item(153,1).Value = "05-08-2021"
$ws.Cells.Item(153,2).Value = 2.7
$ws.Cells.Item(153,3).Value = 3.06
$ws.Cells.Item(153,4).Value = 3.36
$ws.Cells.Item(153,5).Value = 3.61
$ws.Cells.Item(153,6).Value = 4.4
$ws.Cells.Item(153,7).Value = -1.34
$ws.Cells.Item(153,8).Value = -0.66
$ws.Cells.Item(153,9).Value = -0.23
$ws.Cells.Item(153,10).Value = 0.1
$ws.Cells.Item(153,11).Value = 0.36
$ws.Cells.Item(153,12).Value = 1.18
$ws.Cells.Item(153,13).Value = 1.37

$ws.Cells.Item(154,1).NumberFormat = "@"
$ws.Cells.Item(154,1).Value = "06-08-2021"
$ws.Cells.Item(154,2).Value = 2.84
$ws.Cells.Item(154,3).Value = 3.17
$ws.Cells.Item(154,4).Value = 3.46
$ws.Cells.Item(154,5).Value = 3.72
$ws.Cells.Item(154,6).Value = 4.46
$ws.Cells.Item(154,7).Value = -1.55
$ws.Cells.Item(154,8).Value = -0.66
$ws.Cells.Item(154,9).Value = -0.21
$ws.Cells.Item(154,10).Value = 0.15
$ws.Cells.Item(154,11).Value = 0.41
$ws.Cells.Item(154,12).Value = 1.21
$ws.Cells.Item(154,13).Value = 1.42

$ws.Cells.Item(155,1).NumberFormat = "@"
$ws.Cells.Item(155,1).Value = "09-08-2021"
$ws.Cells.Item(155,2).Value = 2.82
$ws.Cells.Item(155,3).Value = 3.15
$ws.Cells.Item(155,4).Value = 3.45
$ws.Cells.Item(155,5).Value = 3.71
$ws.Cells.Item(155,6).Value = 4.46
$ws.Cells.Item(155,7).Value = -1.52
$ws.Cells.Item(155,8).Value = -0.68
$ws.Cells.Item(155,9).Value = -0.23
$ws.Cells.Item(155,10).Value = 0.14
$ws.Cells.Item(155,11).Value = 0.41
$ws.Cells.Item(155,12).Value = 1.2
$ws.Cells.Item(155,13).Value = 1.42

$ws.Cells.Item(156,1).NumberFormat = "@"
$ws.Cells.Item(156,1).Value = "10-08-2021"
$ws.Cells.Item(156,2).Value = 2.8
$ws.Cells.Item(156,3).Value = 3.13
$ws.Cells.Item(156,4).Value = 3.44
$ws.Cells.Item(156,5).Value = 3.7
$ws.Cells.Item(156,6).Value = 4.45
$ws.Cells.Item(156,7).Value = -1.52
$ws.Cells.Item(156,8).Value = -0.7
$ws.Cells.Item(156,9).Value = -0.25
$ws.Cells.Item(156,10).Value = 0.14
$ws.Cells.Item(156,11).Value = 0.41
$ws.Cells.Item(156,12).Value = 1.2
$ws.Cells.Item(156,13).Value = 1.41

$ws.Cells.Item(157,1).NumberFormat = "@"
$ws.Cells.Item(157,1).Value = "11-08-2021"
$ws.Cells.Item(157,2).Value = 2.9
$ws.Cells.Item(157,3).Value = 3.23
$ws.Cells.Item(157,4).Value = 3.52
$ws.Cells.Item(157,5).Value = 3.77
$ws.Cells.Item(157,6).Value = 4.49
$ws.Cells.Item(157,7).Value = -1.66
$ws.Cells.Item(157,8).Value = -0.77
$ws.Cells.Item(157,9).Value = -0.31
$ws.Cells.Item(157,10).Value = 0.09
$ws.Cells.Item(157,11).Value = 0.36
$ws.Cells.Item(157,12).Value = 1.14
$ws.Cells.Item(157,13).Value = 1.37

$ws.Cells.Item(158,1).NumberFormat = "@"
$ws.Cells.Item(158,1).Value = "12-08-2021"
$ws.Cells.Item(158,2).Value = 2.91
$ws.Cells.Item(158,3).Value = 3.25
$ws.Cells.Item(158,4).Value = 3.54
$ws.Cells.Item(158,5).Value = 3.79
$ws.Cells.Item(158,6).Value = 4.51
$ws.Cells.Item(158,7).Value = -1.86
$ws.Cells.Item(158,8).Value = -0.85
$ws.Cells.Item(158,9).Value = -0.42
$ws.Cells.Item(158,10).Value = 0.03
$ws.Cells.Item(158,11).Value = 0.31
$ws.Cells.Item(158,12).Value = 1.12
$ws.Cells.Item(158,13).Value = 1.35

$ws.Cells.Item(159,1).NumberFormat = "@"
$ws.Cells.Item(159,1).Value = "13-08-2021"
$ws.Cells.Item(159,2).Value = 2.93
$ws.Cells.Item(159,3).Value = 3.27
$ws.Cells.Item(159,4).Value = 3.53
$ws.Cells.Item(159,5).Value = 3.78
$ws.Cells.Item(159,6).Value = 4.48
$ws.Cells.Item(159,7).Value = -1.94
$ws.Cells.Item(159,8).Value = -0.93
$ws.Cells.Item(159,9).Value = -0.49
$ws.Cells.Item(159,10).Value = -0.06
$ws.Cells.Item(159,11).Value = 0.24
$ws.Cells.Item(159,12).Value = 1.07
$ws.Cells.Item(159,13).Value = 1.3

$ws.Cells.Item(160,1).NumberFormat = "@"
$ws.Cells.Item(160,1).Value = "16-08-2021"
$ws.Cells.Item(160,2).Value = 2.9
$ws.Cells.Item(160,3).Value = 3.23
$ws.Cells.Item(160,4).Value = 3.5
$ws.Cells.Item(160,5).Value = 3.74
$ws.Cells.Item(160,6).Value = 4.44
$ws.Cells.Item(160,7).Value = -2.08
$ws.Cells.Item(160,8).Value = -1.1
$ws.Cells.Item(160,9).Value = -0.65
$ws.Cells.Item(160,10).Value = -0.2
$ws.Cells.Item(160,11).Value = 0.11
$ws.Cells.Item(160,12).Value = 0.99
$ws.Cells.Item(160,13).Value = 1.23

$ws.Cells.Item(161,1).NumberFormat = "@"
$ws.Cells.Item(161,1).Value = "17-08-2021"
$ws.Cells.Item(161,2).Value = 2.9
$ws.Cells.Item(161,3).Value = 3.23
$ws.Cells.Item(161,4).Value = 3.49
$ws.Cells.Item(161,5).Value = 3.74
$ws.Cells.Item(161,6).Value = 4.44
$ws.Cells.Item(161,7).Value = -2.02
$ws.Cells.Item(161,8).Value = -1.02
$ws.Cells.Item(161,9).Value = -0.57
$ws.Cells.Item(161,10).Value = -0.12
$ws.Cells.Item(161,11).Value = 0.19
$ws.Cells.Item(161,12).Value = 1.05
$ws.Cells.Item(161,13).Value = 1.27

$ws.Cells.Item(162,1).NumberFormat = "@"
$ws.Cells.Item(162,1).Value = "18-08-2021"
$ws.Cells.Item(162,2).Value = 2.97
$ws.Cells.Item(162,3).Value = 3.31
$ws.Cells.Item(162,4).Value = 3.58
$ws.Cells.Item(162,5).Value = 3.81
$ws.Cells.Item(162,6).Value = 4.5
$ws.Cells.Item(162,7).Value = -1.89
$ws.Cells.Item(162,8).Value = -0.88
$ws.Cells.Item(162,9).Value = -0.44
$ws.Cells.Item(162,10).Value = 0.03
$ws.Cells.Item(162,11).Value = 0.33
$ws.Cells.Item(162,12).Value = 1.18
$ws.Cells.Item(162,13).Value = 1.39

$ws.Cells.Item(163,1).NumberFormat = "@"
$ws.Cells.Item(163,1).Value = "19-08-2021"
$ws.Cells.Item(163,2).Value = 3.01
$ws.Cells.Item(163,3).Value = 3.34
$ws.Cells.Item(163,4).Value = 3.61
$ws.Cells.Item(163,5).Value = 3.85
$ws.Cells.Item(163,6).Value = 4.51
$ws.Cells.Item(163,7).Value = -1.88
$ws.Cells.Item(163,8).Value = -0.86
$ws.Cells.Item(163,9).Value = -0.42
$ws.Cells.Item(163,10).Value = 0.04
$ws.Cells.Item(163,11).Value = 0.35
$ws.Cells.Item(163,12).Value = 1.17
$ws.Cells.Item(163,13).Value = 1.37

$ws.Cells.Item(164,1).NumberFormat = "@"
$ws.Cells.Item(164,1).Value = "20-08-2021"
$ws.Cells.Item(164,2).Value = 3.06
$ws.Cells.Item(164,3).Value = 3.39
$ws.Cells.Item(164,4).Value = 3.65
$ws.Cells.Item(164,5).Value = 3.89
$ws.Cells.Item(164,6).Value = 4.54
$ws.Cells.Item(164,7).Value = -1.87
$ws.Cells.Item(164,8).Value = -0.86
$ws.Cells.Item(164,9).Value = -0.42
$ws.Cells.Item(164,10).Value = 0.04
$ws.Cells.Item(164,11).Value = 0.35
$ws.Cells.Item(164,12).Value = 1.18
$ws.Cells.Item(164,13).Value = 1.4

$ws.Cells.Item(165,1).NumberFormat = "@"
$ws.Cells.Item(165,1).Value = "23-08-2021"
$ws.Cells.Item(165,2).Value = 3.06
$ws.Cells.Item(165,3).Value = 3.39
$ws.Cells.Item(165,4).Value = 3.65
$ws.Cells.Item(165,5).Value = 3.89
$ws.Cells.Item(165,6).Value = 4.53
$ws.Cells.Item(165,7).Value = -1.88
$ws.Cells.Item(165,8).Value = -0.89
$ws.Cells.Item(165,9).Value = -0.44
$ws.Cells.Item(165,10).Value = 0.02
$ws.Cells.Item(165,11).Value = 0.32
$ws.Cells.Item(165,12).Value = 1.14
$ws.Cells.Item(165,13).Value = 1.35

$ws.Cells.Item(166,1).NumberFormat = "@"
$ws.Cells.Item(166,1).Value = "24-08-2021"
$ws.Cells.Item(166,2).Value = 3.04
$ws.Cells.Item(166,3).Value = 3.37
$ws.Cells.Item(166,4).Value = 3.62
$ws.Cells.Item(166,5).Value = 3.85
$ws.Cells.Item(166,6).Value = 4.49
$ws.Cells.Item(166,7).Value = -1.91
$ws.Cells.Item(166,8).Value = -0.92
$ws.Cells.Item(166,9).Value = -0.46
$ws.Cells.Item(166,10).Value = 0
$ws.Cells.Item(166,11).Value = 0.31
$ws.Cells.Item(166,12).Value = 1.1
$ws.Cells.Item(166,13).Value = 1.31

$ws.Cells.Item(167,1).NumberFormat = "@"
$ws.Cells.Item(167,1).Value = "25-08-2021"
$ws.Cells.Item(167,2).Value = 3.01
$ws.Cells.Item(167,3).Value = 3.33
$ws.Cells.Item(167,4).Value = 3.58
$ws.Cells.Item(167,5).Value = 3.81
$ws.Cells.Item(167,6).Value = 4.44
$ws.Cells.Item(167,7).Value = -1.91
$ws.Cells.Item(167,8).Value = -0.96
$ws.Cells.Item(167,9).Value = -0.51
$ws.Cells.Item(167,10).Value = -0.04
$ws.Cells.Item(167,11).Value = 0.27
$ws.Cells.Item(167,12).Value = 1.06
$ws.Cells.Item(167,13).Value = 1.26

$ws.Cells.Item(168,1).NumberFormat = "@"
$ws.Cells.Item(168,1).Value = "26-08-2021"
$ws.Cells.Item(168,2).Value = 3.04
$ws.Cells.Item(168,3).Value = 3.37
$ws.Cells.Item(168,4).Value = 3.62
$ws.Cells.Item(168,5).Value = 3.85
$ws.Cells.Item(168,6).Value = 4.45
$ws.Cells.Item(168,7).Value = -1.9
$ws.Cells.Item(168,8).Value = -0.96
$ws.Cells.Item(168,9).Value = -0.51
$ws.Cells.Item(168,10).Value = -0.04
$ws.Cells.Item(168,11).Value = 0.27
$ws.Cells.Item(168,12).Value = 1.08
$ws.Cells.Item(168,13).Value = 1.28

$ws.Cells.Item(169,1).NumberFormat = "@"
$ws.Cells.Item(169,1).Value = "27-08-2021"
$ws.Cells.Item(169,2).Value = 3.07
$ws.Cells.Item(169,3).Value = 3.39
$ws.Cells.Item(169,4).Value = 3.64
$ws.Cells.Item(169,5).Value = 3.87
$ws.Cells.Item(169,6).Value = 4.46
$ws.Cells.Item(169,7).Value = -1.95
$ws.Cells.Item(169,8).Value = -0.96
$ws.Cells.Item(169,9).Value = -0.5
$ws.Cells.Item(169,10).Value = -0.01
$ws.Cells.Item(169,11).Value = 0.31
$ws.Cells.Item(169,12).Value = 1.11
$ws.Cells.Item(169,13).Value = 1.3

$ws.Cells.Item(170,1).NumberFormat = "@"
$ws.Cells.Item(170,1).Value = "30-08-2021"
$ws.Cells.Item(170,2).Value = 3.09
$ws.Cells.Item(170,3).Value = 3.41
$ws.Cells.Item(170,4).Value = 3.65
$ws.Cells.Item(170,5).Value = 3.87
$ws.Cells.Item(170,6).Value = 4.48
$ws.Cells.Item(170,7).Value = -1.95
$ws.Cells.Item(170,8).Value = -0.96
$ws.Cells.Item(170,9).Value = -0.5
$ws.Cells.Item(170,10).Value = -0.01
$ws.Cells.Item(170,11).Value = 0.32
$ws.Cells.Item(170,12).Value = 1.12
$ws.Cells.Item(170,13).Value = 1.32

$ws.Cells.Item(171,1).NumberFormat = "@"
$ws.Cells.Item(171,1).Value = "31-08-2021"
$ws.Cells.Item(171,2).Value = 3.12
$ws.Cells.Item(171,3).Value = 3.45
$ws.Cells.Item(171,4).Value = 3.69
$ws.Cells.Item(171,5).Value = 3.89
$ws.Cells.Item(171,6).Value = 4.53
$ws.Cells.Item(171,7).Value = -1.93
$ws.Cells.Item(171,8).Value = -0.84
$ws.Cells.Item(171,9).Value = -0.43
$ws.Cells.Item(171,10).Value = 0.06
$ws.Cells.Item(171,11).Value = 0.39
$ws.Cells.Item(171,12).Value = 1.18
$ws.Cells.Item(171,13).Value = 1.39

$ws.Cells.Item(172,1).NumberFormat = "@"
$ws.Cells.Item(172,1).Value = "01-09-2021"
$ws.Cells.Item(172,2).Value = 3.7
$ws.Cells.Item(172,3).Value = 3.97
$ws.Cells.Item(172,4).Value = 4.13
$ws.Cells.Item(172,5).Value = 4.29
$ws.Cells.Item(172,6).Value = 4.79
$ws.Cells.Item(172,7).Value = -1.25
$ws.Cells.Item(172,8).Value = -0.36
$ws.Cells.Item(172,9).Value = 0.05
$ws.Cells.Item(172,10).Value = 0.52
$ws.Cells.Item(172,11).Value = 0.8
$ws.Cells.Item(172,12).Value = 1.49
$ws.Cells.Item(172,13).Value = 1.64

$ws.Cells.Item(173,1).NumberFormat = "@"
$ws.Cells.Item(173,1).Value = "02-09-2021"
$ws.Cells.Item(173,2).Value = 3.77
$ws.Cells.Item(173,3).Value = 4.01
$ws.Cells.Item(173,4).Value = 4.17
$ws.Cells.Item(173,5).Value = 4.34
$ws.Cells.Item(173,6).Value = 4.78
$ws.Cells.Item(173,7).Value = -1.09
$ws.Cells.Item(173,8).Value = -0.25
$ws.Cells.Item(173,9).Value = 0.13
$ws.Cells.Item(173,10).Value = 0.58
$ws.Cells.Item(173,11).Value = 0.85
$ws.Cells.Item(173,12).Value = 1.46
$ws.Cells.Item(173,13).Value = 1.62

$ws.Cells.Item(174,1).NumberFormat = "@"
$ws.Cells.Item(174,1).Value = "03-09-2021"
$ws.Cells.Item(174,2).Value = 3.8
$ws.Cells.Item(174,3).Value = 4.05
$ws.Cells.Item(174,4).Value = 4.2
$ws.Cells.Item(174,5).Value = 4.36
$ws.Cells.Item(174,6).Value = 4.82
$ws.Cells.Item(174,7).Value = -1.26
$ws.Cells.Item(174,8).Value = -0.27
$ws.Cells.Item(174,9).Value = 0.11
$ws.Cells.Item(174,10).Value = 0.58
$ws.Cells.Item(174,11).Value = 0.85
$ws.Cells.Item(174,12).Value = 1.44
$ws.Cells.Item(174,13).Value = 1.6

$ws.Range("A150:M174").ClearFormats()
